$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark attendance "O" for week 6 (column G) on the two rows that were
# missing it (row 4 - 소프트웨어공학, row 8 - e-비지니스전략), matching
# the existing "O" already present in column F for those rows.
$ws.Range("G4").Value = "O"
$ws.Range("G8").Value = "O"

# Update the active selection to reflect where the user ended up (G9).
$ws.Range("G9").Select()
